$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 18 as a copy of the existing row 17 (preserves formatting,
# including the date number format used in column D), then update the values
# of each row: row 17 gets the new weekly data, row 18 keeps the original
# (previous week) data that used to live in row 17.

$ws.Rows.Item(17).Copy()
$ws.Rows.Item(18).Insert()

# Row 17: updated values (new week)
$ws.Cells.Item(17, 4).Value = 44706   # D17 Fecha
$ws.Cells.Item(17, 10).Value = 200    # J17 Volumen
$ws.Cells.Item(17, 11).Value = 9000   # K17 Precio minimo
$ws.Cells.Item(17, 12).Value = 9000   # L17 Precio maximo
$ws.Cells.Item(17, 13).Value = 9000   # M17 Precio promedio ponderado
$ws.Cells.Item(17, 16).Value = 250    # P17 Precio $/Kg

# Row 18: original values (previous week), same as row 17 had before the edit
$ws.Cells.Item(18, 1).Value = 5
$ws.Cells.Item(18, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(18, 3).Value = "Maule"
$ws.Cells.Item(18, 4).Value = 44376
$ws.Cells.Item(18, 5).Value = 7
$ws.Cells.Item(18, 6).Value = 100112040
$ws.Cells.Item(18, 7).Value = "Cilantro"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 150
$ws.Cells.Item(18, 11).Value = 6500
$ws.Cells.Item(18, 12).Value = 6500
$ws.Cells.Item(18, 13).Value = 6500
$ws.Cells.Item(18, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(18, 15).Value = "Región Metropolitana"
$ws.Cells.Item(18, 16).Value = 181
$ws.Cells.Item(18, 17).Value = 36
$ws.Cells.Item(18, 18).Value = "Hortaliza"
